$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the PATH value in B2 (the local absolute path is no longer stored in the sheet)
$ws.Range("B2").ClearContents()

# Move the selection to B15, matching the saved cursor position
$ws.Range("B15").Select()
